$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Inline citation text -> pandoc/quarto citation keys
# ------------------------------------------------------------------
$d.Content.Find.Execute("Marrero et al. (2019)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "@marrero2019", 2)

$d.Content.Find.Execute("(Marrero et al. 2019)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[@marrero2019]", 2)

# ------------------------------------------------------------------
# 2. Remove the auto-generated "References" heading + bibliography
#    entry that Quarto/Pandoc appended at the end of the document
#    (the rendered reference list is no longer wanted in the body).
# ------------------------------------------------------------------
$paras = $d.Paragraphs
$refHeadingIndex = -1

for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    $styleName = $p.Range.Style.NameLocal
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($styleName -eq "Heading 2" -and $text -eq "References") {
        $refHeadingIndex = $i
        break
    }
}

if ($refHeadingIndex -ge 1) {
    $startPos = $paras.Item($refHeadingIndex).Range.Start
    $endPos = $d.Content.End
    $killRange = $d.Range($startPos, $endPos)
    $killRange.Delete()
}
